# New weekly price observation was inserted as a new record (row 44),
# pushing all the subsequent "Vega Modelo de Temuco - Espinaca" rows
# down by one. We reproduce that by inserting a blank row at row 44 and
# filling it in with the new observation's data; Excel shifts every
# row below it (44-118 -> 45-119) automatically, which matches the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44 (shifts old rows 44:118 down to 45:119)
$ws.Rows("44:44").Insert()

# Populate the new row 44 with the new weekly record
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = "2022-01-19"
$ws.Range("D44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100112012
$ws.Range("G44").Value = "Espinaca"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 20
$ws.Range("K44").Value = 12000
$ws.Range("L44").Value = 12000
$ws.Range("M44").Value = 12000
$ws.Range("N44").Value = "$/docena de atados"
$ws.Range("O44").Value = "Región de La Araucanía"
$ws.Range("P44").Value = 4000
$ws.Range("Q44").Value = 3
$ws.Range("R44").Value = "Hortaliza"
